$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProblems = @(
    "Insertion in BST",
    "Search given key in BST",
    "Deletion from BST",
    "Construct balanced BST from given keys",
    "Determine if given Binary tree is a BST or not",
    "check if given keys represent same bst without building bst",
    "find inorder predecessor for given key in a bst",
    "find lowest common ansector of two nodes in bst",
    "kth smallest and kth largest element in BST",
    "fllor and ciel in bst",
    "find optimal cost to construct bst",
    "tree to bst while maintiang original strucute",
    "remove nodes from bst that have key outside a given a range",
    "find a pair with given sum in bst",
    "inorder sucessor of bst",
    "fix a binary tree that is onle one swap away from becoming a bst",
    "update every key in bst to contain sum of all greater keys"
)

$startRow = 217
for ($i = 0; $i -lt $newProblems.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newProblems[$i]
}

$ws.Range("B5").Select()
